# Update the cardinal/ordinal number words from French to Spanish on both
# sheets, then restore the selections/active-sheet state recorded by Excel
# when the edit was made.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "1-10"
$ws2 = $wb.Worksheets.Item(2)   # "11-100"

# --- Sheet "1-10": French -> Spanish ---------------------------------
$ws1.Range("C2").Value = "primero"
$ws1.Range("E2").Value = "seis"
$ws1.Range("F2").Value = "sexto"
$ws1.Range("B3").Value = "dos"
$ws1.Range("C3").Value = "segundo"
$ws1.Range("E3").Value = "siete"
$ws1.Range("F3").Value = "séptimo"
$ws1.Range("B4").Value = "tres"
$ws1.Range("C4").Value = "tercero"
$ws1.Range("E4").Value = "ocho"
$ws1.Range("F4").Value = "octavo"
$ws1.Range("B5").Value = "cuarto"
$ws1.Range("C5").Value = "cuatro"
$ws1.Range("E5").Value = "nueve"
$ws1.Range("F5").Value = "noveno"
$ws1.Range("B6").Value = "cinco"
$ws1.Range("C6").Value = "quinto"
$ws1.Range("E6").Value = "diez"
$ws1.Range("F6").Value = "décimo"

# --- Sheet "11-100": French -> Spanish --------------------------------
$ws2.Range("B2").Value = "once"
$ws2.Range("C2").Value = "undécimo"
$ws2.Range("E2").Value = "cuarenta"
$ws2.Range("F2").Value = "cuadra-gésimo"
$ws2.Range("B3").Value = "doce"
$ws2.Range("C3").Value = "duodécimo"
$ws2.Range("E3").Value = "cincuenta"
$ws2.Range("F3").Value = "quincua-gésimo"
$ws2.Range("B4").Value = "trece"
$ws2.Range("C4").Value = "decimotercero"
$ws2.Range("E4").Value = "sesenta"
$ws2.Range("F4").Value = "sexa-gésimo"
$ws2.Range("B5").Value = "catorce"
$ws2.Range("C5").Value = "decimocuarto"
$ws2.Range("E5").Value = "setenta"
$ws2.Range("F5").Value = "septua-gésimo"
$ws2.Range("B6").Value = "quince"
$ws2.Range("C6").Value = "decimoquinto"
$ws2.Range("E6").Value = "ochenta"
$ws2.Range("F6").Value = "octo-gésimo"
$ws2.Range("B7").Value = "veinte"
$ws2.Range("C7").Value = "vi-gésimo"
$ws2.Range("E7").Value = "noventa"
$ws2.Range("F7").Value = "nona-gésimo"
$ws2.Range("B8").Value = "treinta"
$ws2.Range("C8").Value = "tri-gésimo"
$ws2.Range("E8").Value = "ciento"
$ws2.Range("F8").Value = "centésimo"

# --- Selections / active sheet -----------------------------------------
# Sheet "1-10" ends up not tab-selected, with F2:F6 highlighted.
$ws1.Range("F2:F6").Select()

# Sheet "11-100" becomes the tab-selected / active sheet, with F7 the
# active cell (workbook.xml bookViews activeTab -> 1 follows from this).
$ws2.Range("F7").Select()
